$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "73.233.54"
$ws.Range("E2").Value = "  +2.96%  "

$ws.Range("D3").Value = "4.004.64"
$ws.Range("E3").Value = "  +1.04%  "

Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.14%  "

Set-TextValue "D5" "597.29"
$ws.Range("E5").Value = "  +10.78%  "

Set-TextValue "D6" "160.32"
$ws.Range("E6").Value = "  +7.45%  "

$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  -0.07%  "

Set-TextValue "D9" "0.751"
$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("E10").Value = "  +2.36%  "

Set-TextValue "D11" "54.14"
$ws.Range("E11").Value = "  -3.31%  "

Set-TextValue "D12" "0.0000321"
$ws.Range("E12").Value = "  +1.49%  "

Set-TextValue "D13" "11.03"
$ws.Range("E13").Value = "  +3.37%  "

$ws.Range("D14").Value = "4.640.77"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").Value = "4.010.01"
$ws.Range("E15").Value = "  +1.04%  "

Set-TextValue "D16" "1.27"
$ws.Range("E16").Value = "  +8.87%  "

Set-TextValue "D17" "14.18"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D20").Value = "72.903.81"
$ws.Range("E20").Value = "  +2.49%  "

Set-TextValue "D21" "435.79"
$ws.Range("E21").Value = "  +2.10%  "

Set-TextValue "D22" "4.79"
$ws.Range("E22").Value = "  +12.56%  "

$ws.Range("E23").Value = "  -0.91%  "

$ws.Range("E24").Value = "  -4.24%  "

Set-TextValue "D25" "14.26"
$ws.Range("E25").Value = "  -1.50%  "

Set-TextValue "D26" "4.40"
$ws.Range("E26").Value = "  +15.79%  "

$ws.Range("E27").Value = "  -1.32%  "

Set-TextValue "D28" "5.96"
$ws.Range("E28").Value = "  +1.06%  "

Set-TextValue "D29" "10.51"
$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("E30").Value = "  -0.14%  "

Set-TextValue "D31" "7.83"
$ws.Range("E31").Value = "  +0.34%  "

Set-TextValue "D32" "13.81"
$ws.Range("E32").Value = "  +2.95%  "

Set-TextValue "D33" "0.131"
$ws.Range("E33").Value = "  -0.35%  "

Set-TextValue "D34" "48.25"
$ws.Range("E34").Value = "  -4.69%  "

Set-TextValue "D35" "671.83"
$ws.Range("E35").Value = "  -1.78%  "

Set-TextValue "D36" "70.85"
$ws.Range("E36").Value = "  +8.44%  "

$ws.Range("D37").Value = "0.0₃0915"
$ws.Range("E37").Value = "  +11.61%  "

$ws.Range("E38").Value = "  -0.50%  "

Set-TextValue "D39" "0.999"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  -2.67%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D42" "3.33"
$ws.Range("E42").Value = "  +4.17%  "

$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("E44").Value = "  +0.12%  "

Set-TextValue "D45" "10.64"
$ws.Range("E45").Value = "  +9.02%  "

$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "3.40"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D48" "2.61"
$ws.Range("E48").Value = "  -2.75%  "

$ws.Range("D49").Value = "2.876.28"
$ws.Range("E49").Value = "  +8.53%  "

$ws.Range("E50").Value = "  +1.63%  "

$ws.Range("E51").Value = "  +4.39%  "
